$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (single decimal point) need to be
# forced to Text format first, otherwise Excel auto-converts the literal
# string into a floating point number (losing the exact printed text,
# e.g. "576.06" -> 576.05999999999995, or "46.00" -> 46).
$textForceCells = @('D4', 'D5', 'D6', 'D9', 'D10', 'D14', 'D18', 'D19', 'D21', 'D23', 'D24', 'D27', 'D28', 'D29', 'D32', 'D35', 'D37', 'D38', 'D40', 'D42', 'D44', 'D45', 'D47', 'D49')
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.875.24'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '3.073.90'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '576.06'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = '167.47'
$ws.Range('E6').Value = '  -2.79%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.072.14'
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('D9').Value = '0.513'
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('D10').Value = '6.38'
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('E12').Value = '  -1.99%  '
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('D14').Value = '35.99'
$ws.Range('E14').Value = '  -3.51%  '
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.781.95'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '3.579.32'
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').Value = '7.01'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('D19').Value = '16.85'
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('D20').Value = '3.067.28'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').Value = '484.32'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('E22').Value = '  -3.54%  '
$ws.Range('D23').Value = '7.68'
$ws.Range('E23').Value = '  -4.13%  '
$ws.Range('D24').Value = '82.57'
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('E25').Value = '  -4.93%  '
$ws.Range('E26').Value = '  -2.99%  '
$ws.Range('D27').Value = '10.13'
$ws.Range('E27').Value = '  +2.57%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = '7.83'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('E30').Value = '  -6.78%  '
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('D32').Value = '27.69'
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('E33').Value = '  -2.07%  '
$ws.Range('D34').Value = '0.0₃0899'
$ws.Range('E34').Value = '  -3.21%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  -3.12%  '
$ws.Range('D37').Value = '0.953'
$ws.Range('E37').Value = '  -2.11%  '
$ws.Range('D38').Value = '46.00'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').Value = '1.99'
$ws.Range('E40').Value = '  -4.41%  '
$ws.Range('E41').Value = '  -2.14%  '
$ws.Range('D42').Value = '8.32'
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('D43').Value = '2.767.87'
$ws.Range('E43').Value = '  -1.13%  '
$ws.Range('D44').Value = '369.02'
$ws.Range('E44').Value = '  -2.14%  '
$ws.Range('D45').Value = '135.87'
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('E46').Value = '  -3.06%  '
$ws.Range('D47').Value = '2.46'
$ws.Range('E47').Value = '  -2.78%  '
$ws.Range('D49').Value = '24.42'
$ws.Range('E49').Value = '  -1.39%  '
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('E51').Value = '  -1.78%  '

# Restore the default (unstyled) cell style now that the text is stored,
# so these cells do not carry a lingering explicit Text number format.
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
